$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=-0.09914768942982544; C=0.6760636891099804; D=1.058759192261776;  E=1.02896024814459;   F=1.052237025233689;  G=19},
    @{Row=3;  B=-0.1143386297721907;  C=0.6889840097079422; D=0.9692665404913205; E=0.9845133521142924; F=1.006200714633619;  G=18},
    @{Row=4;  B=-0.03231471279853988; C=0.7695157101427289; D=0.9504010021523444; E=0.9748851225412891; F=1.004336377456358;  G=17},
    @{Row=5;  B=0.05671622629984657;  C=0.6321910596607474; D=0.8518300650585396; E=0.9229464042177854; F=0.951413452211165;  G=16},
    @{Row=6;  B=0.085928964333323;    C=0.5612671955956002; D=0.5873845025919633; E=0.7664101399328973; F=0.7883078955770884; G=15},
    @{Row=7;  B=0.1382138394267656;   C=0.7476276879241753; D=1.097975515794379;  E=1.047843268716452;  F=1.077897322974905;  G=14},
    @{Row=8;  B=0.1215206328360093;   C=0.7819526689445239; D=1.205421270947194;  E=1.097916786895616;  F=1.135726716235405;  G=13},
    @{Row=9;  B=0.2028461735751207;   C=0.8472735835465873; D=1.158512388347609;  E=1.076342133500128;  F=1.10405820904526;   G=12},
    @{Row=10; B=0.239907996146195;    C=0.6651045147971374; D=0.5928596743010739; E=0.769973814035954;  F=0.7673553351966808; G=11},
    @{Row=11; B=0.2225839694609408;   C=0.8306889768957216; D=1.245987321698922;  E=1.116238021973325;  F=1.152988242321246;  G=10}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$wb.Save()
